$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Duplicate the current row 4 (Mehmet / Özyürek / 5055555555 / Akademisyen) down to row 5,
#    carrying only its values (not its formatting) - it becomes the new last row.
$ws.Range("A4:D4").Copy()
$ws.Range("A5:D5").PasteSpecial(-4104)

# 2. Overwrite row 4 with the newly added person (Yunus Şen).
$ws.Range("A4").Value = "Yunus"
$ws.Range("B4").Value = "Şen"
$ws.Range("C4").Value = 5305555555
$ws.Range("D4").Value = "Öğrenci"

# 3. Update the phone numbers for the first two existing data rows.
$ws.Range("C2").Value = 5305555555
$ws.Range("C3").Value = 5305555555

# 4. Give column A (the name column) on every data row the same bordered / centered
#    formatting the header row already uses.
$ws.Range("A1").Copy()
$ws.Range("A2:A5").PasteSpecial(-4122)

# 5. The header no longer needs to be bold - match the (now shared) look of column A.
$ws.Range("A2:A5").Font.Bold = $false
$ws.Range("A1:D1").Font.Bold = $false

# 6. Cosmetic selection change.
$ws.Range("E4").Select()
